# Refresh "Pais" (countries) COVID stats table and timestamp banner.
# Column layout: A=Pais, B=Casos totales, C=Nuevos casos, D=Casos activos,
#                E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes.
# The data is kept sorted descending by column B ("Casos totales"), so
# refreshed totals shuffle a handful of rows up/down; those rows also get
# their country name (column A) rewritten to match the new sort position.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1: "last updated" banner
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 21 de Julio de 2020 a las 20:18"

# Row 4
$ws.Cells.Item(4, 2).Value = 3989363
$ws.Cells.Item(4, 3).Value = 27934
$ws.Cells.Item(4, 4).Value = 1862360
$ws.Cells.Item(4, 5).Value = 1982599
$ws.Cells.Item(4, 7).Value = 570
$ws.Cells.Item(4, 8).Value = 144404
# Row 6
$ws.Cells.Item(6, 2).Value = 1192151
$ws.Cells.Item(6, 3).Value = 37234
$ws.Cells.Item(6, 4).Value = 752312
$ws.Cells.Item(6, 5).Value = 411070
$ws.Cells.Item(6, 7).Value = 670
$ws.Cells.Item(6, 8).Value = 28769
# Row 18
$ws.Cells.Item(18, 2).Value = 221500
$ws.Cells.Item(18, 3).Value = 928
$ws.Cells.Item(18, 4).Value = 204011
$ws.Cells.Item(18, 5).Value = 11963
$ws.Cells.Item(18, 7).Value = 18
$ws.Cells.Item(18, 8).Value = 5526
# Row 21
$ws.Cells.Item(21, 2).Value = 203597
$ws.Cells.Item(21, 3).Value = 110
$ws.Cells.Item(21, 5).Value = 6322
# Row 24
$ws.Cells.Item(24, 2).Value = 111508
$ws.Cells.Item(24, 3).Value = 384
$ws.Cells.Item(24, 4).Value = 97674
$ws.Cells.Item(24, 5).Value = 4974
$ws.Cells.Item(24, 7).Value = 2
$ws.Cells.Item(24, 8).Value = 8860
# Row 59
$ws.Cells.Item(59, 2).Value = 25802
$ws.Cells.Item(59, 3).Value = 36
$ws.Cells.Item(59, 5).Value = 685
# Row 61
$ws.Cells.Item(61, 2).Value = 24278
$ws.Cells.Item(61, 3).Value = 587
$ws.Cells.Item(61, 4).Value = 16646
$ws.Cells.Item(61, 5).Value = 6532
$ws.Cells.Item(61, 7).Value = 13
$ws.Cells.Item(61, 8).Value = 1100
# Row 67
$ws.Cells.Item(67, 2).Value = 17742
$ws.Cells.Item(67, 3).Value = 180
$ws.Cells.Item(67, 4).Value = 15389
$ws.Cells.Item(67, 5).Value = 2073
$ws.Cells.Item(67, 7).Value = 4
$ws.Cells.Item(67, 8).Value = 280
# Row 70
$ws.Cells.Item(70, 1).Value = "Kenia"
$ws.Cells.Item(70, 2).Value = 14168
$ws.Cells.Item(70, 3).Value = 397
$ws.Cells.Item(70, 4).Value = 6258
$ws.Cells.Item(70, 5).Value = 7660
$ws.Cells.Item(70, 7).Value = 12
$ws.Cells.Item(70, 8).Value = 250
# Row 71
$ws.Cells.Item(71, 1).Value = "Chequia"
$ws.Cells.Item(71, 2).Value = 14160
$ws.Cells.Item(71, 3).Value = 62
$ws.Cells.Item(71, 4).Value = 8899
$ws.Cells.Item(71, 5).Value = 4901
$ws.Cells.Item(71, 7).Value = 1
$ws.Cells.Item(71, 8).Value = 360
# Row 72
$ws.Cells.Item(72, 1).Value = "Corea del Sur"
$ws.Cells.Item(72, 2).Value = 13816
$ws.Cells.Item(72, 3).Value = 45
$ws.Cells.Item(72, 4).Value = 12643
$ws.Cells.Item(72, 5).Value = 877
$ws.Cells.Item(72, 8).Value = 296
# Row 110
$ws.Cells.Item(110, 2).Value = 3044
$ws.Cells.Item(110, 3).Value = 45
$ws.Cells.Item(110, 4).Value = 2397
$ws.Cells.Item(110, 5).Value = 632
# Row 111
$ws.Cells.Item(111, 2).Value = 2980
$ws.Cells.Item(111, 3).Value = 75
$ws.Cells.Item(111, 4).Value = 1577
$ws.Cells.Item(111, 5).Value = 1362
# Row 119
$ws.Cells.Item(119, 2).Value = 2107
$ws.Cells.Item(119, 3).Value = 36
$ws.Cells.Item(119, 4).Value = 1100
$ws.Cells.Item(119, 5).Value = 986
# Row 120
$ws.Cells.Item(120, 1).Value = "Libia"
$ws.Cells.Item(120, 2).Value = 2088
$ws.Cells.Item(120, 3).Value = 108
$ws.Cells.Item(120, 4).Value = 479
$ws.Cells.Item(120, 5).Value = 1559
$ws.Cells.Item(120, 7).Value = 1
$ws.Cells.Item(120, 8).Value = 50
# Row 121
$ws.Cells.Item(121, 1).Value = "Estonia"
$ws.Cells.Item(121, 2).Value = 2022
$ws.Cells.Item(121, 3).Value = 1
$ws.Cells.Item(121, 4).Value = 1912
$ws.Cells.Item(121, 5).Value = 41
$ws.Cells.Item(121, 8).Value = 69
# Row 122
$ws.Cells.Item(122, 1).Value = "Eslovaquia"
$ws.Cells.Item(122, 2).Value = 2021
$ws.Cells.Item(122, 3).Value = 41
$ws.Cells.Item(122, 4).Value = 1538
$ws.Cells.Item(122, 5).Value = 455
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 28
# Row 123
$ws.Cells.Item(123, 1).Value = "Hong Kong"
$ws.Cells.Item(123, 2).Value = 2020
$ws.Cells.Item(123, 3).Value = 61
$ws.Cells.Item(123, 4).Value = 1324
$ws.Cells.Item(123, 5).Value = 682
$ws.Cells.Item(123, 7).Value = 2
$ws.Cells.Item(123, 8).Value = 14
# Row 128
$ws.Cells.Item(128, 2).Value = 1894
$ws.Cells.Item(128, 3).Value = 68
$ws.Cells.Item(128, 4).Value = 855
$ws.Cells.Item(128, 5).Value = 1015
$ws.Cells.Item(128, 7).Value = 1
$ws.Cells.Item(128, 8).Value = 24
# Row 129
$ws.Cells.Item(129, 1).Value = "Sierra Leona"
$ws.Cells.Item(129, 2).Value = 1727
$ws.Cells.Item(129, 3).Value = 16
$ws.Cells.Item(129, 4).Value = 1273
$ws.Cells.Item(129, 5).Value = 388
$ws.Cells.Item(129, 8).Value = 66
# Row 130
$ws.Cells.Item(130, 1).Value = "Zimbabue"
$ws.Cells.Item(130, 2).Value = 1713
$ws.Cells.Item(130, 4).Value = 472
$ws.Cells.Item(130, 5).Value = 1215
$ws.Cells.Item(130, 8).Value = 26
# Row 148
$ws.Cells.Item(148, 4).Value = 803
$ws.Cells.Item(148, 5).Value = 11
# Row 157
$ws.Cells.Item(157, 2).Value = 639
$ws.Cells.Item(157, 3).Value = 8
$ws.Cells.Item(157, 5).Value = 164
# Row 165
$ws.Cells.Item(165, 4).Value = 280
$ws.Cells.Item(165, 5).Value = 55
# Row 184
$ws.Cells.Item(184, 2).Value = 111
$ws.Cells.Item(184, 3).Value = 2
$ws.Cells.Item(184, 5).Value = 7
# Row 200
$ws.Cells.Item(200, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(200, 3).Value = 8
$ws.Cells.Item(200, 4).Value = 8
$ws.Cells.Item(200, 5).Value = 18
$ws.Cells.Item(200, 8).Value = 1
# Row 201
$ws.Cells.Item(201, 1).Value = "Fiyi"
$ws.Cells.Item(201, 2).Value = 27
$ws.Cells.Item(201, 4).Value = 18
$ws.Cells.Item(201, 5).Value = 9
# Row 202
$ws.Cells.Item(202, 1).Value = "Timor Oriental"
$ws.Cells.Item(202, 2).Value = 24
$ws.Cells.Item(202, 4).Value = 24
$ws.Cells.Item(202, 5).Value = 0
# Row 203
$ws.Cells.Item(203, 1).Value = "Santa Lucia"
$ws.Cells.Item(203, 4).Value = 19
$ws.Cells.Item(203, 5).Value = 4
# Row 204
$ws.Cells.Item(204, 1).Value = "Granada"
$ws.Cells.Item(204, 2).Value = 23
$ws.Cells.Item(204, 4).Value = 23
$ws.Cells.Item(204, 5).Value = 0
# Row 205
$ws.Cells.Item(205, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(205, 2).Value = 22
$ws.Cells.Item(205, 4).Value = 21
$ws.Cells.Item(205, 5).Value = 1
$ws.Cells.Item(205, 8).Value = 0
